$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: set a cell value while forcing text storage so that
# numeric-looking strings (e.g. "576.17") are NOT reinterpreted as numbers,
# and restore the original cell style afterwards so no formatting changes leak in.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$updates = @(
    @{ Cell = 'D2'; Value = '63.755.01'; ForceText = $true }
    @{ Cell = 'E2'; Value = '  +1.77%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.475.55'; ForceText = $true }
    @{ Cell = 'E3'; Value = '  +1.56%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.01%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '576.17'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +1.58%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '148.60'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +2.32%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.06%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.545'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  +1.97%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '2.475.98'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +1.56%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.112'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +1.12%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  +1.16%  '; ForceText = $false }
    @{ Cell = 'B12'; Value = 'Toncoin'; ForceText = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false }
    @{ Cell = 'D12'; Value = '5.28'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'B13'; Value = 'Cardano'; ForceText = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.359'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +1.39%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '27.32'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +1.64%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.0000182'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -1.27%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '2.924.97'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +1.65%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '63.551.36'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  +1.84%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '2.447.69'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  +0.29%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '11.52'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +2.53%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '7.40'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +6.71%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '330.71'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +2.11%  '; ForceText = $false }
    @{ Cell = 'E22'; Value = '  +1.32%  '; ForceText = $false }
    @{ Cell = 'E23'; Value = '  +20.21%  '; ForceText = $false }
    @{ Cell = 'E24'; Value = '  +0.10%  '; ForceText = $false }
    @{ Cell = 'E25'; Value = '  -1.44%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '631.40'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +11.77%  '; ForceText = $false }
    @{ Cell = 'B27'; Value = 'Aptos'; ForceText = $false }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false }
    @{ Cell = 'D27'; Value = '9.08'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +3.92%  '; ForceText = $false }
    @{ Cell = 'B28'; Value = 'PEPE'; ForceText = $false }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; ForceText = $false }
    @{ Cell = 'D28'; Value = '0.0000106'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  +3.94%  '; ForceText = $false }
    @{ Cell = 'B29'; Value = 'WrappedeETH'; ForceText = $false }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; ForceText = $false }
    @{ Cell = 'D29'; Value = '2.603.44'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +1.95%  '; ForceText = $false }
    @{ Cell = 'B30'; Value = 'Fetch.AI'; ForceText = $false }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; ForceText = $false }
    @{ Cell = 'D30'; Value = '1.53'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +5.24%  '; ForceText = $false }
    @{ Cell = 'E31'; Value = '  +0.13%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '8.40'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +0.14%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '0.144'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -2.30%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '1.91'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +2.47%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '5.28'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +8.29%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '1.55'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +0.37%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  +0.17%  '; ForceText = $false }
    @{ Cell = 'E38'; Value = '  +0.27%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '5.49'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +1.09%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '18.86'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +0.44%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '2.75'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +13.74%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  +0.30%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '147.37'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -0.73%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -0.22%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '150.79'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +1.45%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '3.79'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +2.97%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '21.50'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +4.71%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '0.0544'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +1.15%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '0.607'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +1.20%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.0236'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +2.14%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.0921'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -0.65%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        Set-TextValue $range $u.Value
    } else {
        $range.Value = $u.Value
    }
}

